# Sync workbook data: insert two new quotations (and their items), and
# fix the productService_type on an existing item row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quotations")
$ws2 = $wb.Worksheets.Item("items")

# ---------------------------------------------------------------------
# Sheet "quotations": insert a new row at 5 ("Orçamento teste", pending)
# ---------------------------------------------------------------------
$ws1.Rows("5:5").Insert()
$ws1.Cells.Item(5,1).Value  = "NzZhZGQwNGUtMmU5ZS00NjM5LTkzYzgtNWQxOTc3OGJlYWRlOjU3MDE2"
$ws1.Cells.Item(5,2).Value  = "BKHV0XAFEG"
$ws1.Cells.Item(5,3).Value  = "Orçamento teste"
$ws1.Cells.Item(5,4).Value  = "Reclamaçao do cliente :nao atinge temperatura"
$ws1.Cells.Item(5,5).Value  = "Link do Relatório: https://app.fieldcontrol.com.br/relacionamento-cliente/#/relatorio/657fd93e-b9fd-457f-9f3a-db4c162fa0d7"
$ws1.Cells.Item(5,6).Value  = $false
$ws1.Cells.Item(5,7).Value  = "'1535"
$ws1.Cells.Item(5,8).Value  = "'1535"
$ws1.Cells.Item(5,9).Value  = "Pendente"
$ws1.Cells.Item(5,10).Value = "2025-09-04T16:01:06.048Z"
$ws1.Cells.Item(5,13).Value = "Adriana Vieira Masini"
$ws1.Cells.Item(5,15).Value = "2025-08-28T16:05:56.490Z"
$ws1.Cells.Item(5,17).Value = "percentage"
$ws1.Cells.Item(5,18).Value = "'0"
$ws1.Cells.Item(5,19).Value = "'0"
$ws1.Cells.Item(5,20).Value = "NDIzOTk2OTo1NzAxNg=="
$ws1.Cells.Item(5,21).Value = "pending"

# ---------------------------------------------------------------------
# Sheet "quotations": insert a new row at 7 ("teste", approved)
# (row numbers below are the ones that exist *after* the prior insert)
# ---------------------------------------------------------------------
$ws1.Rows("7:7").Insert()
$ws1.Cells.Item(7,1).Value  = "OTcwY2YxNjktMmI1NS00ZWVkLWI3YmMtYzQxOGQzODUyMmY0OjU3MDE2"
$ws1.Cells.Item(7,2).Value  = "S1SWKWAKXG"
$ws1.Cells.Item(7,3).Value  = "teste"
$ws1.Cells.Item(7,4).Value  = "Reclamaçao do cliente: nao atinge temperatura, abaixo segue link"
$ws1.Cells.Item(7,6).Value  = $false
$ws1.Cells.Item(7,7).Value  = "'1535"
$ws1.Cells.Item(7,8).Value  = "'1535"
$ws1.Cells.Item(7,9).Value  = "Aprovada"
$ws1.Cells.Item(7,10).Value = "2025-09-04T16:07:52.834Z"
$ws1.Cells.Item(7,12).Value = "rffccfc"
$ws1.Cells.Item(7,13).Value = "Adriana Vieira Masini"
$ws1.Cells.Item(7,14).Value = "accounts/57016/quotations/970cf169-2b55-4eed-b7bc-c418d38522f4/signatures/7a22a671-9310-4602-ba1e-f1006682889f.png"
$ws1.Cells.Item(7,15).Value = "2025-08-28T16:10:11.398Z"
$ws1.Cells.Item(7,17).Value = "percentage"
$ws1.Cells.Item(7,18).Value = "'0"
$ws1.Cells.Item(7,19).Value = "'0"
$ws1.Cells.Item(7,20).Value = "NDIzOTk2OTo1NzAxNg=="
$ws1.Cells.Item(7,21).Value = "approved"

# ---------------------------------------------------------------------
# Sheet "items": row 3 productService_type fix (product -> service)
# ---------------------------------------------------------------------
$ws2.Cells.Item(3,9).Value = "service"

# ---------------------------------------------------------------------
# Sheet "items": insert two new rows at 6-7 for the "Orçamento teste" quotation
# ---------------------------------------------------------------------
$ws2.Rows("6:7").Insert()

$ws2.Cells.Item(6,1).Value  = "NTIzYmQ3M2MtMzQwOC00ZTM3LWJhYzctZDM2YzZjYjg0NjY5OjU3MDE2"
$ws2.Cells.Item(6,2).Value  = 1
$ws2.Cells.Item(6,3).Value  = 680
$ws2.Cells.Item(6,4).Value  = "WRVT.00021 REALIZADO SERVIÇO LIMPEZA E CARGA DE GAS  R`$680,00"
$ws2.Cells.Item(6,5).Value  = 3
$ws2.Cells.Item(6,6).Value  = "NzZhZGQwNGUtMmU5ZS00NjM5LTkzYzgtNWQxOTc3OGJlYWRlOjU3MDE2"
$ws2.Cells.Item(6,7).Value  = "NWIwZWFlMmYtYjVkOC00NTU0LTkzZmYtZGM2ZGIwM2E1ZmEwOjU3MDE2"
$ws2.Cells.Item(6,8).Value  = 680
$ws2.Cells.Item(6,9).Value  = "service"
$ws2.Cells.Item(6,10).Value = "NzZhZGQwNGUtMmU5ZS00NjM5LTkzYzgtNWQxOTc3OGJlYWRlOjU3MDE2"

$ws2.Cells.Item(7,1).Value  = "NjU4NTlmOGEtZWM4NC00NGE1LTliMmQtNzcwMzRkOTI0ZDMwOjU3MDE2"
$ws2.Cells.Item(7,2).Value  = 1
$ws2.Cells.Item(7,3).Value  = 855
$ws2.Cells.Item(7,4).Value  = "WRVT.00020  REALIZADO RECUPERAÇAO DA ESTAÇAO MICRO MOTOR  E REALIZADO LIMPEZA  NO SISTEMA E CARGA DE GAS"
$ws2.Cells.Item(7,5).Value  = 3
$ws2.Cells.Item(7,6).Value  = "NzZhZGQwNGUtMmU5ZS00NjM5LTkzYzgtNWQxOTc3OGJlYWRlOjU3MDE2"
$ws2.Cells.Item(7,7).Value  = "MWY3MGI1MWUtZWEwMC00YWEyLTgzZTItNDgwYzc2NzE1OTJkOjU3MDE2"
$ws2.Cells.Item(7,8).Value  = 855
$ws2.Cells.Item(7,9).Value  = "service"
$ws2.Cells.Item(7,10).Value = "NzZhZGQwNGUtMmU5ZS00NjM5LTkzYzgtNWQxOTc3OGJlYWRlOjU3MDE2"

# ---------------------------------------------------------------------
# Sheet "items": insert two new rows at 9-10 for the "teste" quotation
# (row numbers below are the ones that exist *after* the prior insert)
# ---------------------------------------------------------------------
$ws2.Rows("9:10").Insert()

$ws2.Cells.Item(9,1).Value   = "NTMwYmZhMDQtZGU2Ny00ZmFlLWE2NmItNmM5YzMxODdiZTgwOjU3MDE2"
$ws2.Cells.Item(9,2).Value   = 1
$ws2.Cells.Item(9,3).Value   = 680
$ws2.Cells.Item(9,4).Value   = "WRVT.00021 REALIZADO SERVIÇO LIMPEZA E CARGA DE GAS  R`$680,00"
$ws2.Cells.Item(9,5).Value   = 3
$ws2.Cells.Item(9,6).Value   = "OTcwY2YxNjktMmI1NS00ZWVkLWI3YmMtYzQxOGQzODUyMmY0OjU3MDE2"
$ws2.Cells.Item(9,7).Value   = "NWIwZWFlMmYtYjVkOC00NTU0LTkzZmYtZGM2ZGIwM2E1ZmEwOjU3MDE2"
$ws2.Cells.Item(9,8).Value   = 680
$ws2.Cells.Item(9,9).Value   = "service"
$ws2.Cells.Item(9,10).Value  = "OTcwY2YxNjktMmI1NS00ZWVkLWI3YmMtYzQxOGQzODUyMmY0OjU3MDE2"

$ws2.Cells.Item(10,1).Value  = "ZTFiNzY1Y2MtNTVhMS00NDU4LTliMmMtMTU3MzM1ZjU5ZTA4OjU3MDE2"
$ws2.Cells.Item(10,2).Value  = 1
$ws2.Cells.Item(10,3).Value  = 855
$ws2.Cells.Item(10,4).Value  = "WRVT.00020  REALIZADO RECUPERAÇAO DA ESTAÇAO MICRO MOTOR  E REALIZADO LIMPEZA  NO SISTEMA E CARGA DE GAS"
$ws2.Cells.Item(10,5).Value  = 3
$ws2.Cells.Item(10,6).Value  = "OTcwY2YxNjktMmI1NS00ZWVkLWI3YmMtYzQxOGQzODUyMmY0OjU3MDE2"
$ws2.Cells.Item(10,7).Value  = "MWY3MGI1MWUtZWEwMC00YWEyLTgzZTItNDgwYzc2NzE1OTJkOjU3MDE2"
$ws2.Cells.Item(10,8).Value  = 855
$ws2.Cells.Item(10,9).Value  = "service"
$ws2.Cells.Item(10,10).Value = "OTcwY2YxNjktMmI1NS00ZWVkLWI3YmMtYzQxOGQzODUyMmY0OjU3MDE2"
